# Corrected Thingiverse and GitHub links
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New combined label/URL text for the 3D-printed "Make" rows (12, 13, 14)
$label = "Thingiverse`nAlternate: GitHub"
$url = "https://www.thingiverse.com/thing:4790412`nalternate: https://github.com/MirageC79/HextrudORT/tree/main/files/CARRIAGE/HD12/STL"

$ws.Range("J12").Value = $label
$ws.Range("K12").Value = $url

$ws.Range("J13").Value = $label
$ws.Range("K13").Value = $url

$ws.Range("J14").Value = $label
$ws.Range("K14").Value = $url

# Widen column K to fit new content
$ws.Columns.Item(11).ColumnWidth = 41.7265625

# Update the view: scroll/zoom and selection to I16
$ws.Application.ActiveWindow.Zoom = 85
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("I16").Select()
